# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" (Overview sheet) and the
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" columns
# (zh-cn and de-de sheets) for the da5e691d...md row, reflecting the
# newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to da5e691d-f6b3-4ecb-a132-7b63c4c63144.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G3").Value = "2016-09-06 08:39:29"

# zh-cn sheet: row 3 corresponds to da5e691d-f6b3-4ecb-a132-7b63c4c63144.md
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$wsZhCn.Range("H3").Value = "2016-09-06 08:39:17"
$wsZhCn.Range("K3").Value = "2016-09-06 08:40:33"

# de-de sheet: row 3 corresponds to da5e691d-f6b3-4ecb-a132-7b63c4c63144.md
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$wsDeDe.Range("H3").Value = "2016-09-06 08:39:29"
$wsDeDe.Range("K3").Value = "2016-09-06 08:40:52"
